$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the old "Comprehension scores" column (C),
# pushing it and "What I did" (D) one column to the right (D, E). ---
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 13.5

# The new column's header ("Time spent that week") belongs next to "Week",
# so swap it into B1 and push the old B1 ("Total time") header into the
# freshly inserted C1 cell.
$oldB1 = $ws.Range("B1").Text
$ws.Range("C1").Value = $oldB1
$ws.Range("B1").Value = "Time spent that week"

# Fill in the running-total formulas for the new column C, rows 2-14
# (existing weeks), using the same [h]:mm:ss elapsed-time format as B.
$ws.Range("C2:C14").NumberFormat = "[h]:mm:ss"
$ws.Range("C2").Formula = "=SUM(B2)+1.2708333333"
$ws.Range("C3").Formula = "=SUM(B2:B3)+1.2708333333"
$ws.Range("C4").Formula = "=SUM(B2:B4)+1.2708333333"
$ws.Range("C5").Formula = "=SUM(B2:B5)+1.2708333333"
$ws.Range("C6").Formula = "=SUM(B2:B6)+1.2708333333"
$ws.Range("C7").Formula = "=SUM(B2:B7)+1.2708333333"
$ws.Range("C8").Formula = "=SUM(B2:B8)+1.2708333333"
$ws.Range("C9").Formula = "=SUM(B2:B9)+1.2708333333"
$ws.Range("C10").Formula = "=SUM(B2:B10)+1.2708333333"
$ws.Range("C11").Formula = "=SUM(B2:B11)+1.2708333333"
$ws.Range("C12").Formula = "=SUM(B2:B12)+1.2708333333"
$ws.Range("C13").Formula = "=SUM(B2:B13)+1.2708333333"
$ws.Range("C14").Formula = "=SUM(B2:B14)+1.2708333333"

# A trailing space was appended to the week-13 comprehension-scores text (D14).
$d14 = $ws.Range("D14")
$d14.Value = $d14.Text + " "

# --- Add week 14 as a new row. ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 0.88201388888888888
$ws.Range("B15").NumberFormat = "h:mm:ss"
$ws.Range("C15").NumberFormat = "[h]:mm:ss"
$ws.Range("C15").Formula = "=SUM(B2:B15)+1.2708333333"
$ws.Range("D15").Value = "Lord of the Rings (Audiovisual, English, Familiar):36; La casa de las flores (Subtitled, Spanish, New):33; Oscuro deseo (Audiovisual, Spanish, Re-watch):35;  ¿Quién mató a Sara? (Audiovisual, Spanish, Re-watch):38;"

# Update the active selection to match the new last row.
$ws.Range("C15").Select()
